# "updated user management page"
#
# Adds a new "UserDetails" worksheet (after the existing "logins" sheet)
# that holds a single user's profile (username/password/employee
# name/status/role), and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the "logins" sheet so the new sheet inherits the same
# sheetView/sheetFormatPr/namespace scaffolding used throughout this
# workbook, then place it right after "logins" and rename it.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "UserDetails"

# Populate the new user-details table.
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Andreson"
$ws2.Range("A2").Value = "Password"
$ws2.Range("B2").Value = "anil123"
$ws2.Range("A3").Value = "EmployeeName"
$ws2.Range("B3").Value = "Kevin  Mathews"
$ws2.Range("A4").Value = "Status"
$ws2.Range("B4").Value = "Enabled"
$ws2.Range("A5").Value = "UserRole"
$ws2.Range("B5").Value = "ESS"

# Widen the label column so the longer labels aren't clipped.
$ws2.Columns.Item(1).ColumnWidth = 26.6667

# Match the new sheet's page setup (inches: 0.75/0.75/1/1/0.5/0.5).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Leave the cursor on B1 and make "UserDetails" the active/selected tab.
$ws2.Range("B1").Select() | Out-Null
$ws2.Activate() | Out-Null
